# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.217.96'
$ws.Range("E2").Value = '  +4.20%  '
$ws.Range("D3").Value = '2.324.88'
$ws.Range("E3").Value = '  +1.41%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '522.00'
$ws.Range("E5").Value = '  +3.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.15'
$ws.Range("E6").Value = '  +4.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("E8").Value = '  +1.47%  '
$ws.Range("D9").Value = '2.349.85'
$ws.Range("E9").Value = '  +1.63%  '
$ws.Range("E10").Value = '  +5.79%  '
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.31'
$ws.Range("E12").Value = '  +4.56%  '
$ws.Range("E13").Value = '  +0.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.97'
$ws.Range("E14").Value = '  +0.72%  '
$ws.Range("D15").Value = '2.753.18'
$ws.Range("E15").Value = '  +1.90%  '
$ws.Range("D16").Value = '57.026.56'
$ws.Range("E16").Value = '  +3.81%  '
$ws.Range("E17").Value = '  +2.23%  '
$ws.Range("D18").Value = '2.344.10'
$ws.Range("E18").Value = '  +2.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.50'
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.23'
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.05'
$ws.Range("E21").Value = '  +4.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.62'
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.88'
$ws.Range("E24").Value = '  +1.02%  '
$ws.Range("E25").Value = '  +8.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.988'
$ws.Range("E26").Value = '  -0.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.97'
$ws.Range("E27").Value = '  +6.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.30'
$ws.Range("E28").Value = '  +14.12%  '
$ws.Range("E29").Value = '  +5.49%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.73'
$ws.Range("E30").Value = '  +5.90%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.44'
$ws.Range("E31").Value = '  -2.67%  '
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.36'
$ws.Range("E33").Value = '  +1.77%  '
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("E35").Value = '  +0.32%  '
$ws.Range("E36").Value = '  +1.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.929'
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("E38").Value = '  +3.80%  '
$ws.Range("E39").Value = '  +7.50%  '
$ws.Range("E40").Value = '  +3.04%  '
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.60'
$ws.Range("E42").Value = '  +4.71%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '138.47'
$ws.Range("E43").Value = '  +3.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.23'
$ws.Range("E44").Value = '  +1.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '277.93'
$ws.Range("E45").Value = '  +7.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0934'
$ws.Range("E46").Value = '  +2.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0507'
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.564'
$ws.Range("E48").Value = '  +2.56%  '
$ws.Range("E49").Value = '  +3.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.84'
$ws.Range("E50").Value = '  +7.82%  '
$ws.Range("E51").Value = '  +0.66%  '
